$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "purchase"
$ws.Range("A9").Value = "restore from clean install"
$ws.Range("A10").Value = "thumbs"
$ws.Range("A11").Value = "splashes"

$ws.Range("B4").Select()
